$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados..." timestamp footer
$ws.Range("A1").Value = "Datos actualizados a 12 de Julio de 2020 a las 18:08"

# Country name reorderings caused by updated case counts changing rank order
$ws.Range("A81").Value = "Republica de Macedonia"
$ws.Range("A82").Value = "Consejo Danes para los Refugiados"
$ws.Range("A96").Value = "Luxemburgo"
$ws.Range("A97").Value = "Madagascar"
$ws.Range("A113").Value = "Libano"
$ws.Range("A114").Value = "Malaui"
$ws.Range("A134").Value = "Montenegro"
$ws.Range("A135").Value = "Jordania"
$ws.Range("A136").Value = "Letonia"
$ws.Range("A169").Value = "Lesoto"
$ws.Range("A170").Value = "Eritrea"
$ws.Range("A171").Value = "Mongolia"
$ws.Range("A172").Value = "Islas Caimanes"
$ws.Range("A173").Value = "Burundi"
$ws.Range("A174").Value = "Guadalupe"
$ws.Range("A175").Value = "Islas Feroe"

# Updated numeric statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B4").Value = 3381730
$ws.Range("C4").Value = 26084
$ws.Range("D4").Value = 1502215
$ws.Range("E4").Value = 1741937
$ws.Range("G4").Value = 175
$ws.Range("H4").Value = 137578
$ws.Range("B5").Value = 1846249
$ws.Range("C5").Value = 5437
$ws.Range("E5").Value = 561153
$ws.Range("G5").Value = 92
$ws.Range("H5").Value = 71584
$ws.Range("B6").Value = 871499
$ws.Range("C6").Value = 21141
$ws.Range("D6").Value = 546379
$ws.Range("E6").Value = 302042
$ws.Range("G6").Value = 391
$ws.Range("H6").Value = 23078
$ws.Range("B12").Value = 289603
$ws.Range("C12").Value = 650
$ws.Range("G12").Value = 21
$ws.Range("H12").Value = 44819
$ws.Range("B16").Value = 243061
$ws.Range("C16").Value = 234
$ws.Range("D16").Value = 194928
$ws.Range("E16").Value = 13179
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = 34954
$ws.Range("B19").Value = 199903
$ws.Range("C19").Value = 91
$ws.Range("E19").Value = 6269
$ws.Range("B23").Value = 107589
$ws.Range("C23").Value = 242
$ws.Range("D23").Value = 71467
$ws.Range("E23").Value = 27339
$ws.Range("G23").Value = 10
$ws.Range("H23").Value = 8783
$ws.Range("B42").Value = 46512
$ws.Range("C42").Value = 291
$ws.Range("D42").Value = 30907
$ws.Range("E42").Value = 13945
$ws.Range("G42").Value = 6
$ws.Range("H42").Value = 1660
$ws.Range("D43").Value = 42285
$ws.Range("E43").Value = 3650
$ws.Range("B68").Value = 13148
$ws.Range("C68").Value = 33
$ws.Range("D68").Value = 8229
$ws.Range("E68").Value = 4567
$ws.Range("B81").Value = 8111
$ws.Range("C81").Value = 136
$ws.Range("D81").Value = 4203
$ws.Range("E81").Value = 3526
$ws.Range("G81").Value = 6
$ws.Range("H81").Value = 382
$ws.Range("B82").Value = 8033
$ws.Range("C82").Value = 62
$ws.Range("D82").Value = 3615
$ws.Range("E82").Value = 4229
$ws.Range("H82").Value = 189
$ws.Range("B96").Value = 4925
$ws.Range("C96").Value = 83
$ws.Range("D96").Value = 4086
$ws.Range("E96").Value = 728
$ws.Range("H96").Value = 111
$ws.Range("B97").Value = 4867
$ws.Range("C97").Value = 289
$ws.Range("D97").Value = 2378
$ws.Range("E97").Value = 2454
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 35
$ws.Range("B100").Value = 3803
$ws.Range("C100").Value = 31
$ws.Range("E100").Value = 2236
$ws.Range("B110").Value = 2612
$ws.Range("C110").Value = 101
$ws.Range("E110").Value = 620
$ws.Range("B113").Value = 2334
$ws.Range("C113").Value = 166
$ws.Range("D113").Value = 1420
$ws.Range("E113").Value = 878
$ws.Range("H113").Value = 36
$ws.Range("B114").Value = 2261
$ws.Range("D114").Value = 517
$ws.Range("E114").Value = 1711
$ws.Range("H114").Value = 33
$ws.Range("B134").Value = 1221
$ws.Range("C134").Value = 57
$ws.Range("D134").Value = 325
$ws.Range("E134").Value = 873
$ws.Range("H134").Value = 23
$ws.Range("B135").Value = 1179
$ws.Range("C135").Value = 3
$ws.Range("D135").Value = 997
$ws.Range("E135").Value = 172
$ws.Range("H135").Value = 10
$ws.Range("B136").Value = 1173
$ws.Range("D136").Value = 1019
$ws.Range("E136").Value = 124
$ws.Range("H136").Value = 30
$ws.Range("B141").Value = 1021
$ws.Range("C141").Value = 7
$ws.Range("E141").Value = 163
$ws.Range("B169").Value = 233
$ws.Range("C169").Value = 49
$ws.Range("D169").Value = 32
$ws.Range("E169").Value = 199
$ws.Range("G169").Value = 1
$ws.Range("H169").Value = 2
$ws.Range("B170").Value = 232
$ws.Range("C170").Value = 0
$ws.Range("D170").Value = 107
$ws.Range("E170").Value = 125
$ws.Range("B171").Value = 230
$ws.Range("C171").Value = 3
$ws.Range("D171").Value = 202
$ws.Range("E171").Value = 28
$ws.Range("H171").Value = 0
$ws.Range("B172").Value = 201
$ws.Range("D172").Value = 197
$ws.Range("E172").Value = 3
$ws.Range("B173").Value = 191
$ws.Range("D173").Value = 118
$ws.Range("E173").Value = 72
$ws.Range("H173").Value = 1
$ws.Range("B174").Value = 190
$ws.Range("D174").Value = 157
$ws.Range("E174").Value = 19
$ws.Range("H174").Value = 14
$ws.Range("B175").Value = 188
$ws.Range("D175").Value = 188
$ws.Range("E175").Value = 0
$ws.Range("H175").Value = 0
